{"js": "// Office.js (Word JavaScript API) script implementing:\n//  1. Turn the plain-text stackexchange URL in the \"Acceleration discussion - \"\n//     bullet into a real hyperlink (same visible text, Hyperlink style).\n//  2. Append a brand-new References bullet:\n//       \"Parabolic equation = \" + hyperlink to the Parabolas.html page.\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// Part 1: convert the existing stackexchange URL text into a hyperlink\n// ---------------------------------------------------------------------\nconst url1 = \"https://physics.stackexchange.com/questions/513405/acceleration-time-graph-for-a-falling-object\";\n\nconst found1 = body.search(url1, { matchCase: true });\nfound1.load(\"items\");\nawait context.sync();\n\nif (found1.items.length > 0) {\n  found1.items[0].hyperlink = url1;\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// Part 2: add a new list bullet \"Parabolic equation = <hyperlink>\"\n// ---------------------------------------------------------------------\nconst url2 = \"https://www.csun.edu/~ayk38384/notes/mod11/Parabolas.html\";\nconst prefix2 = \"Parabolic equation = \";\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n// Inherits the ListParagraph style + numPr (same numbered/bulleted list)\n// from the paragraph it follows, matching the other References bullets.\nconst newParagraph = lastParagraph.insertParagraph(prefix2 + url2, Word.InsertLocation.after);\nawait context.sync();\n\nconst found2 = body.search(url2, { matchCase: true });\nfound2.load(\"items\");\nawait context.sync();\n\nif (found2.items.length > 0) {\n  found2.items[0].hyperlink = url2;\n  await context.sync();\n}\n", "ps1": "# Word COM interop script implementing:\n#  1. Turn the plain-text URL in the \"Acceleration discussion - \" bullet into\n#     a real hyperlink (same visible text, Hyperlink character style applied).\n#  2. Append a brand-new References bullet:\n#       \"Parabolic equation = \" + hyperlink to the Parabolas.html page.\n#\n# Note: this runtime's Hyperlinks.Add always inserts the new hyperlink run\n# at the *start* of the paragraph that owns the supplied range, no matter\n# where inside the paragraph that range actually points. So for both edits\n# we: (a) make the paragraph contain *only* the URL text, (b) clear that\n# text down to a collapsed range, (c) call Hyperlinks.Add on the now-empty\n# paragraph (so \"start of paragraph\" == the right spot), then (d)\n# InsertBefore() any label text that must precede the link.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# Part 1: convert the existing stackexchange URL text run into a hyperlink\n# ---------------------------------------------------------------------\n$url1 = \"https://physics.stackexchange.com/questions/513405/acceleration-time-graph-for-a-falling-object\"\n$prefix1 = \"Acceleration discussion - \"\n\n$findRange = $d.Content\n$findRange.Find.ClearFormatting()\n$found = $findRange.Find.Execute($url1)\n\nif ($found) {\n    $para = $findRange.Paragraphs(1)\n    $pStart = $para.Range.Start\n    $pEndNoMark = $para.Range.End - 1\n\n    # Clear the whole paragraph's text (run text only, keep paragraph mark)\n    $fullRange = $d.Range($pStart, $pEndNoMark)\n    $fullRange.Text = \"\"\n\n    # Add the hyperlink - it lands at the (now sole) start of the paragraph\n    $emptyRange = $d.Range($pStart, $pStart)\n    $null = $d.Hyperlinks.Add($emptyRange, $url1)\n\n    # Re-insert the label text that must precede the URL\n    $prefixRange = $d.Range($pStart, $pStart)\n    $prefixRange.InsertBefore($prefix1)\n}\n\n# ---------------------------------------------------------------------\n# Part 2: add a new list bullet \"Parabolic equation = <hyperlink>\"\n# ---------------------------------------------------------------------\n$url2 = \"https://www.csun.edu/~ayk38384/notes/mod11/Parabolas.html\"\n$prefix2 = \"Parabolic equation = \"\n\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs($d.Paragraphs.Count)\n# Put the URL text in as the paragraph's sole content (replaces the\n# placeholder run Word creates for a brand-new empty paragraph).\n$newPara.Range.Text = $url2\n\n$newStart = $newPara.Range.Start\n$newEndNoMark = $newPara.Range.End - 1\n$wholeRange = $d.Range($newStart, $newEndNoMark)\n$wholeRange.Text = \"\"\n\n$emptyRange2 = $d.Range($newStart, $newStart)\n$null = $d.Hyperlinks.Add($emptyRange2, $url2)\n\n$prefixRange2 = $d.Range($newStart, $newStart)\n$prefixRange2.InsertBefore($prefix2)\n"}
